$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("parms")

# Update the N parameter row (row 17) values from 500000 to 100000 across columns B:H
$ws.Range("B17:H17").Value = 100000

# Update the active cell selection to A17
$ws.Range("A17").Select()
